$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.413.56'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.27%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.997.72'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +4.18%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '324.20'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.23%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.10%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5099'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.37%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4131'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.40%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08718'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +5.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.133'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.00%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '43.10'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +2.51%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.61'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +3.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.997.74'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.558'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.01%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.443'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +1.90%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.003'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.16%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '94.20'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001114'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06506'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.88'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +3.90%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.05%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.181'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +4.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '30.467.74'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.90'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +5.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.226'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.69%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.229.87'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +4.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.31'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.22%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.98'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.394'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +4.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '131.36'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.133'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.06%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.060'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.838'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.56%  '
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +10.61%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02518'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +2.98%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.434'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.76%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06592'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +2.39%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +9.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2196'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.42%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.029'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.06%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6614'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +2.86%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.233'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.60'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.76%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6156'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.54%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.200'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.26%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.668'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.266'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +4.25%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.45'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.79%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '80.24'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.65%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06893'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.48%  '
